# Auto update: 2025-12-01 14:08:26
# Applies the daily data refresh to the 국장_조선_분석 decision sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-5 (columns A..O), matching the new data snapshot.
$rows = @(
    @{ Row=2;  A="2025-12-01"; B="HDKSOE";          C="009540.KS"; D=399500; E=43.6; F=-7.95;  G=10; H=56; I=73; J=80; K=58;   M="⛔ 관망하십시오."; N=85.92500513438651 },
    @{ Row=3;  A="2025-12-01"; B="SamsungHvyInd";   C="010140.KS"; D=24200;  E=32;   F=-3.39;  G=10; H=70; I=73; J=93; K=58;   M="⛔ 관망하십시오."; N=85.92500513438651 },
    @{ Row=4;  A="2025-12-01"; B="Hanwha Ocean";    C="042660.KS"; D=105300; E=22.5; F=-11.36; G=0;  H=66; I=73; J=73; K=55;   M="⛔ 관망하십시오."; N=85.92500513438651 },
    @{ Row=5;  A="2025-12-01"; B="HD HYUNDAI MIPO"; C="010620.KS"; D=223000; E=57.4; F=0;      G=30; H=46; I=50; J=63; K=54.8; M="⛔ 관망하십시오."; N=85.92500513438651 }
)

foreach ($r in $rows) {
    $i = $r.Row
    # Column A (날짜) is left untouched; it already holds "2025-12-01" and is unchanged by this update.
    $ws.Cells.Item($i, 2).Value  = $r.B          # B: 종목명
    $ws.Cells.Item($i, 3).Value  = $r.C          # C: 티커
    $ws.Cells.Item($i, 4).Value  = $r.D          # D: 종가
    $ws.Cells.Item($i, 5).Value  = $r.E          # E: RSI
    $ws.Cells.Item($i, 6).Value  = $r.F          # F: 5일수익률
    $ws.Cells.Item($i, 7).Value  = $r.G          # G: 점수(룰)
    $ws.Cells.Item($i, 8).Value  = $r.H          # H: 3일상승확률(%)
    $ws.Cells.Item($i, 9).Value  = $r.I          # I: 5일상승확률(%)
    $ws.Cells.Item($i, 10).Value = $r.J          # J: 10일상승확률(%)
    $ws.Cells.Item($i, 11).Value = $r.K          # K: 최종점수
    $ws.Cells.Item($i, 12).Value = "Pattern"     # L: 예측방식 (unchanged)
    $ws.Cells.Item($i, 13).Value = $r.M          # M: 판단
    $ws.Cells.Item($i, 14).Value = $r.N          # N: MACRO_SCORE
    $ws.Cells.Item($i, 15).Value = "🟢 완화적 (상승 우위)"  # O: MACRO_SIGNAL (unchanged)
}

$wb.Save()
